# Streamlines and comments sections
# - Adds two new yearly sheets ("2005" and "2010") as copies of the
#   existing "2000" balance-sheet tab.
# - Makes "2010" the active/selected tab (third tab, index 2).
# - Updates the saved cell selection on each tab.

$wb = $excel.ActiveWorkbook

# The original (and only) worksheet in the workbook.
$ws2000 = $wb.Worksheets.Item("2000")

# Duplicate "2000" to create "2005", placing the copy at the end of the
# worksheet collection, then rename it.
$ws2000.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$wb.Worksheets.Item($wb.Worksheets.Count).Name = "2005"

# Duplicate "2000" again to create "2010", again appended at the end.
$ws2000.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$wb.Worksheets.Item($wb.Worksheets.Count).Name = "2010"

$ws2005 = $wb.Worksheets.Item("2005")
$ws2010 = $wb.Worksheets.Item("2010")

# Restore/update the per-sheet selection state.
$ws2000.Activate()
$ws2000.Range("C37").Select()

$ws2005.Activate()
$ws2005.Range("C37").Select()

# "2010" ends up as the active tab shown when the workbook is opened.
$ws2010.Activate()
$ws2010.Range("I40").Select()
